$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vars = @(
    "_tejgtotfun_f5amb",
    "_tejgfun_f2ct05pgrco",
    "_tejgfun_f5ct05opseg",
    "piagtotfun_f5r07pgrcopc",
    "pimgfun_f5r18ct06cydeppc",
    "_tejgct_r09gstcrpc",
    "piagct_r18gstcppc",
    "pimgfun_f5ct06trans",
    "dfgpimpiatotfun_f5r07pgrco",
    "dfgdevpiagtotfun_f5r18pgrco",
    "tejgtotfun_f5trans",
    "pimgtotfun_f5comepc",
    "tdvgtotfun_f5r18salud",
    "_dfgpimpiatotfun_f5r18opseg",
    "_dfgpimpiafun_f5ct06cydep",
    "_dfgdevpiagct_r18gstcppc",
    "dfgdevpiagfun_f5r18ct05transpc",
    "devppimfun_f5r07ct06cydep",
    "_tejgtotfun_f5r07pgrco",
    "_dfgpimpiatotfun_f5edu",
    "_dfgdevpiagge_r09ct05biserpc",
    "mod_104",
    "tdvgfun_f2ct06salud",
    "_dfgpimpiatotfun_f5r07turi",
    "dfgpimpiafun_f1ct06san",
    "dfgdevpiagtotfun_f5r08cydeppc",
    "dfgdevpiagfun_f5r18ct05pgrco",
    "tdvgfun_f1ct05transpc",
    "_piagfun_f1ct05trans",
    "dfgdevpiagfun_f5r07ct05transpc",
    "_tdvgtotfun_f5r08opseg",
    "_dfgpimpiatotfun_f5amb",
    "_devppimfun_f5ct05salud",
    "_dfgdevpiagtotfun_f5saludpc",
    "_dfgdevpiagfun_f5r07ct05saludpc",
    "pimgfun_f5r18ct06pgrco",
    "_tdvgfun_f3ct06edu",
    "impres_muni_6",
    "dfgpimpiact_r07gstcp",
    "piagfun_f5ct05cydeppc",
    "dfgdevpiagfun_f5ct06prots",
    "per_012",
    "dfgpimpiact_r09gstcr",
    "pimgct_r00gstcp",
    "devppimfun_f5r18ct05cydep",
    "dfgdevpiagtotfun_f5protspc",
    "_pimgfun_f2ct05vivpc",
    "tdvgge_r13ct06acanf",
    "_tejgfun_f5ct05prots",
    "dfgdevpiagfun_f5ct06agropc"
)

$scores = @(
    0.03063949942588806,
    0.02281138114631176,
    0.01305772364139557,
    0.01299246214330196,
    0.01081372890621424,
    0.009939749725162983,
    0.009924954734742641,
    0.009797048754990101,
    0.009533995762467384,
    0.009142960421741009,
    0.00888920109719038,
    0.008313522674143314,
    0.008286050520837307,
    0.007981852628290653,
    0.007303737103939056,
    0.007178127765655518,
    0.007120235823094845,
    0.007120139896869659,
    0.006827239878475666,
    0.006097807548940182,
    0.005942713003605604,
    0.005838421639055014,
    0.005726317409425974,
    0.005561329424381256,
    0.005425289738923311,
    0.005391712766140699,
    0.005294968374073505,
    0.005181117448955774,
    0.005140307359397411,
    0.005004383623600006,
    0.004977420438081026,
    0.004938862752169371,
    0.004903604742139578,
    0.004799069371074438,
    0.00479853805154562,
    0.004752232693135738,
    0.004750834777951241,
    0.004704366438090801,
    0.004691069480031729,
    0.00461532361805439,
    0.004500477574765682,
    0.004479204770177603,
    0.004364656284451485,
    0.004276749677956104,
    0.004265927243977785,
    0.004233721178025007,
    0.004199815448373556,
    0.004188514314591885,
    0.004157342482358217,
    0.004147694911807775
)

for ($i = 0; $i -lt $vars.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $vars[$i]
    $ws.Cells.Item($row, 3).Value = $scores[$i]
}
